$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts existing rows 35:142 down to 36:143,
# matching the target dimension A1:R143.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = [DateTime]"2021-09-14"
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112043
$ws.Range("G35").Value = "Pepino ensalada"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 16000
$ws.Range("L35").Value = 17000
$ws.Range("M35").Value = 16500
$ws.Range("N35").Value = "$/caja 60 unidades"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 275
$ws.Range("Q35").Value = 60
$ws.Range("R35").Value = "Hortaliza"
